# Auto-generated Excel COM-interop script to apply the market-data refresh diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 250000450
$ws.Range("I5").Value = 147.5
$ws.Range("J5").Value = 500000740
$ws.Range("K5").Value = 147.5
$ws.Range("L5").Value = 500000740
$ws.Range("M5").Value = -32.5
$ws.Range("N5").Value = -500000970
$ws.Range("H9").Value = 110450.2
$ws.Range("I9").Value = 183600.33
$ws.Range("J9").Value = 725
$ws.Range("K9").Value = 183600.33
$ws.Range("L9").Value = 725
$ws.Range("M9").Value = -183431.33
$ws.Range("N9").Value = -1063
$ws.Range("H17").Value = 2149
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2149
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6447
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6783
$ws.Range("H33").Value = 299.85715
$ws.Range("I33").Value = 188.75
$ws.Range("K33").Value = 188.75
$ws.Range("M33").Value = 40.25
$ws.Range("H96").Value = 989.63635
$ws.Range("J96").Value = 1185.2
$ws.Range("L96").Value = 3555.6
$ws.Range("N96").Value = -6301.6
$ws.Range("H115").Value = 759.94116
$ws.Range("I115").Value = 759.94116
$ws.Range("K115").Value = 2279.82348
$ws.Range("M115").Value = -712.82348
$ws.Range("H116").Value = 10601.75
$ws.Range("I116").Value = 9973.429
$ws.Range("J116").Value = 15000
$ws.Range("K116").Value = 9973.429
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = -6531.429
$ws.Range("N116").Value = -21884
$ws.Range("H125").Value = 2565.3572
$ws.Range("I125").Value = 1344.75
$ws.Range("K125").Value = 12102.75
$ws.Range("M125").Value = -9642.75
$ws.Range("H134").Value = 31938.934
$ws.Range("J134").Value = 31938.934
$ws.Range("L134").Value = 31938.934
$ws.Range("N134").Value = -42078.934
$ws.Range("H138").Value = 2638.1123
$ws.Range("J138").Value = 2122.6494
$ws.Range("L138").Value = 6367.948199999999
$ws.Range("N138").Value = -16647.9482

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4082.4
$ws.Range("I32").Value = 3804.1355
$ws.Range("K32").Value = 3804.1355
$ws.Range("M32").Value = -3517.1355
$ws.Range("H61").Value = 21514.5
$ws.Range("I61").Value = 3015
$ws.Range("K61").Value = 3015
$ws.Range("M61").Value = -2803
$ws.Range("H74").Value = 2505.8
$ws.Range("I74").Value = 1011
$ws.Range("J74").Value = 3502.3333
$ws.Range("K74").Value = 1011
$ws.Range("L74").Value = 3502.3333
$ws.Range("M74").Value = -137
$ws.Range("N74").Value = -5250.3333
$ws.Range("H77").Value = 2505.8
$ws.Range("I77").Value = 1011
$ws.Range("J77").Value = 3502.3333
$ws.Range("K77").Value = 5055
$ws.Range("L77").Value = 17511.6665
$ws.Range("M77").Value = -687
$ws.Range("N77").Value = -26247.6665
$ws.Range("H97").Value = 2339.8667
$ws.Range("I97").Value = 502.38235
$ws.Range("J97").Value = 8019.364
$ws.Range("K97").Value = 502.38235
$ws.Range("L97").Value = 8019.364
$ws.Range("M97").Value = -6.382349999999974
$ws.Range("N97").Value = -9011.364
$ws.Range("H136").Value = 21514.5
$ws.Range("I136").Value = 3015
$ws.Range("K136").Value = 9045
$ws.Range("M136").Value = -6495

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 7888.5
$ws.Range("I36").Value = 737
$ws.Range("J36").Value = 15040
$ws.Range("K36").Value = 737
$ws.Range("L36").Value = 15040
$ws.Range("M36").Value = -203
$ws.Range("N36").Value = -16108
$ws.Range("H39").Value = 27276.5
$ws.Range("J39").Value = 27276.5
$ws.Range("L39").Value = 27276.5
$ws.Range("N39").Value = -28054.5
$ws.Range("H86").Value = 5952.222
$ws.Range("I86").Value = 5892.6665
$ws.Range("J86").Value = 6071.3335
$ws.Range("K86").Value = 5892.6665
$ws.Range("L86").Value = 6071.3335
$ws.Range("M86").Value = -4769.6665
$ws.Range("N86").Value = -8317.333500000001
$ws.Range("H89").Value = 5952.222
$ws.Range("I89").Value = 5892.6665
$ws.Range("J89").Value = 6071.3335
$ws.Range("K89").Value = 29463.3325
$ws.Range("L89").Value = 30356.6675
$ws.Range("M89").Value = -23847.3325
$ws.Range("N89").Value = -41588.6675

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2352.2856
$ws.Range("I58").Value = 2199.9412
$ws.Range("K58").Value = 2199.9412
$ws.Range("M58").Value = -1996.9412
$ws.Range("H136").Value = 2352.2856
$ws.Range("I136").Value = 2199.9412
$ws.Range("K136").Value = 6599.823600000001
$ws.Range("M136").Value = -4049.823600000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 290.90698
$ws.Range("J2").Value = 781
$ws.Range("L2").Value = 4686
$ws.Range("N2").Value = -4912
$ws.Range("H131").Value = 17778868
$ws.Range("I131").Value = 6945146.5
$ws.Range("J131").Value = 37038816
$ws.Range("K131").Value = 20835439.5
$ws.Range("L131").Value = 111116448
$ws.Range("M131").Value = -20830399.5
$ws.Range("N131").Value = -111126528

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 37984
$ws.Range("J26").Value = 37814.668
$ws.Range("L26").Value = 37814.668
$ws.Range("N26").Value = -38374.668
$ws.Range("H50").Value = 37984
$ws.Range("J50").Value = 37814.668
$ws.Range("L50").Value = 37814.668
$ws.Range("N50").Value = -38810.668
$ws.Range("H70").Value = 11145.167
$ws.Range("I70").Value = 12374.375
$ws.Range("K70").Value = 12374.375
$ws.Range("M70").Value = -12104.375
$ws.Range("H73").Value = 11145.167
$ws.Range("I73").Value = 12374.375
$ws.Range("K73").Value = 12374.375
$ws.Range("M73").Value = -11438.375
$ws.Range("H80").Value = 3696.0476
$ws.Range("J80").Value = 4280.375
$ws.Range("L80").Value = 4280.375
$ws.Range("N80").Value = -6276.375
$ws.Range("H83").Value = 3696.0476
$ws.Range("J83").Value = 4280.375
$ws.Range("L83").Value = 21401.875
$ws.Range("N83").Value = -31385.875
$ws.Range("H122").Value = 5839.7144
$ws.Range("I122").Value = 5571
$ws.Range("J122").Value = 6198
$ws.Range("K122").Value = 16713
$ws.Range("L122").Value = 18594
$ws.Range("M122").Value = -14263
$ws.Range("N122").Value = -23494

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2549.5
$ws.Range("I93").Value = 2433
$ws.Range("J93").Value = 2666
$ws.Range("K93").Value = 2433
$ws.Range("L93").Value = 2666
$ws.Range("M93").Value = -1185
$ws.Range("N93").Value = -5162
$ws.Range("H136").Value = 3326.8965
$ws.Range("I136").Value = 3178.6667
$ws.Range("J136").Value = 3485.7144
$ws.Range("K136").Value = 9536.000100000001
$ws.Range("L136").Value = 10457.1432
$ws.Range("M136").Value = -6986.000100000001
$ws.Range("N136").Value = -15557.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 580.6667
$ws.Range("J113").Value = 193
$ws.Range("L113").Value = 579
$ws.Range("N113").Value = -4919
$ws.Range("H132").Value = 2144.12
$ws.Range("I132").Value = 2051.2
$ws.Range("J132").Value = 2515.8
$ws.Range("K132").Value = 6153.599999999999
$ws.Range("L132").Value = 7547.400000000001
$ws.Range("M132").Value = -3623.599999999999
$ws.Range("N132").Value = -12607.4
$ws.Range("H136").Value = 9067.370000000001
$ws.Range("I136").Value = 7326.1
$ws.Range("K136").Value = 21978.3
$ws.Range("M136").Value = -19428.3
